$d = $word.ActiveDocument

function Append-AfterText {
    param(
        [string]$searchText,
        [string]$appendText
    )
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $insertRange = $d.Range($rng.End, $rng.End)
        $insertRange.InsertAfter($appendText)
    }
}

Append-AfterText "Yang non bulanan berhasil notif Cuma sisa pembayaran blm keluar" " (DONE)"
Append-AfterText "Pembayraan pomg masih muncul template bawaan" " (BISA DICEK ULANG)"
Append-AfterText "Spp juga sama" " (BISA DICEK ULANG)"
Append-AfterText "Spp ta juga sama" " (BISA DICEK ULANG)"
Append-AfterText "Uang kegiatan juga sama" " (BISA DICEK ULANG)"
Append-AfterText "Di login siswa" " (DONE)"
Append-AfterText "Saat edit profile siswa ada bugs" " (DONE)"

Write-Output "done"
